$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "No" column (column A) entirely; this shifts every other
# column (Kode Barang ... Keterangan) one position to the left.
$ws.Range("A1").EntireColumn.Delete()

# Move the active selection to match the edited workbook.
$ws.Range("B9").Select()
